$wb = $excel.ActiveWorkbook

# 1. Update the status text from "Ready for handoff" to "In Translation"
#    across every sheet that contains it.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the "Status" columns that had been widened for the longer text.
#    Target stored width (per the canonical OOXML) is 13.4101845877511
#    characters; this engine quantizes ColumnWidth writes to 1/6-character
#    steps (stored = round(ColumnWidth*6)/6 + 5/6), so 12.5 is the input
#    that lands on the closest reachable stored value (13.333333333333334).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
